$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner (A1)
$ws.Range("A1").Value = "Datos actualizados a 23 de Marzo de 2020 a las 14:16"

# Refresh the full countries table (rows 4-194). Each entry is:
#   row, Pais, Casos totales, Nuevos casos, Casos activos, Recuperados,
#   Casos criticos, Muertes hoy, Muertes
# Countries are now re-sorted by "Casos totales" descending per the
# latest data pull (several countries moved rows, e.g. Chile, Islandia,
# San Marino and Sri Lanka moved up; numeric stats were refreshed too).
$data = @(
    @(4, 'China', 81093, 39, 72703, 5120, 1749, 9, 3270),
    @(5, 'Italia', 59138, 0, 7024, 46638, 3000, 0, 5476),
    @(6, 'Estados Unidos', 35075, 1529, 178, 34439, 795, 39, 458),
    @(7, 'España', 33089, 4321, 3355, 27528, 2355, 434, 2206),
    @(8, 'Alemania', 26220, 1347, 422, 25687, 23, 17, 111),
    @(9, 'Iran', 23049, 1411, 8376, 12861, 0, 127, 1812),
    @(10, 'Francia', 16689, 671, 2200, 13815, 1746, 0, 674),
    @(11, 'Corea del Sur', 8961, 64, 3166, 5684, 59, 7, 111),
    @(12, 'Suiza', 8547, 1073, 131, 8298, 141, 20, 118),
    @(13, 'Reino Unido', 5683, 0, 135, 5259, 20, 8, 289),
    @(14, 'Paises Bajos', 4204, 0, 2, 4023, 354, 0, 179),
    @(15, 'Austria', 3865, 283, 9, 3840, 14, 0, 16),
    @(16, 'Belgica', 3743, 342, 401, 3254, 322, 13, 88),
    @(17, 'Noruega', 2538, 153, 6, 2523, 41, 2, 9),
    @(18, 'Portugal', 2060, 460, 14, 2023, 47, 9, 23),
    @(19, 'Suecia', 1934, 0, 16, 1893, 80, 4, 25),
    @(20, 'Australia', 1717, 108, 88, 1622, 11, 0, 7),
    @(21, 'Brasil', 1620, 74, 2, 1593, 18, 0, 25),
    @(22, 'Malasia', 1518, 212, 159, 1345, 57, 4, 14),
    @(23, 'Canada', 1472, 2, 14, 1438, 1, 0, 20),
    @(24, 'Dinamarca', 1450, 55, 1, 1436, 55, 0, 13),
    @(25, 'Israel', 1238, 167, 37, 1200, 24, 0, 1),
    @(26, 'Turquia', 1236, 0, 0, 1206, 0, 0, 30),
    @(27, 'Chequia', 1165, 45, 6, 1158, 19, 0, 1),
    @(28, 'Japon', 1101, 0, 235, 825, 49, 0, 41),
    @(29, 'Irlanda', 906, 0, 5, 897, 29, 0, 4),
    @(30, 'Pakistan', 824, 48, 13, 805, 0, 1, 6),
    @(31, 'Luxemburgo', 798, 0, 6, 784, 3, 0, 8),
    @(32, 'Ecuador', 789, 0, 3, 772, 2, 0, 14),
    @(33, 'Chile', 746, 114, 8, 737, 7, 0, 1),
    @(34, 'Tailandia', 721, 122, 52, 668, 7, 0, 1),
    @(35, 'Crucero', 712, 0, 567, 137, 15, 0, 8),
    @(36, 'Finlandia', 700, 74, 10, 689, 13, 0, 1),
    @(37, 'Polonia', 684, 50, 13, 663, 3, 1, 8),
    @(38, 'Grecia', 624, 0, 19, 588, 18, 2, 17),
    @(39, 'Islandia', 588, 20, 36, 551, 14, 0, 1),
    @(40, 'Indonesia', 579, 65, 30, 500, 0, 1, 49),
    @(41, 'Rumania', 576, 143, 73, 498, 15, 2, 5),
    @(42, 'Arabia Saudita', 511, 0, 17, 494, 0, 0, 0),
    @(43, 'Singapur', 509, 54, 152, 355, 14, 0, 2),
    @(44, 'Catar', 494, 0, 33, 461, 6, 0, 0),
    @(45, 'Filipinas', 462, 82, 18, 411, 1, 8, 33),
    @(46, 'Rusia', 438, 71, 17, 420, 0, 0, 1),
    @(47, 'India', 425, 29, 24, 393, 0, 1, 8),
    @(48, 'Eslovenia', 414, 0, 0, 411, 12, 1, 3),
    @(49, 'Sudafrica', 402, 128, 2, 400, 0, 0, 0),
    @(50, 'Peru', 363, 0, 1, 357, 5, 0, 5),
    @(51, 'Hong Kong', 356, 38, 100, 252, 4, 0, 4),
    @(52, 'Estonia', 352, 26, 4, 348, 4, 0, 0),
    @(53, 'Barein', 339, 5, 160, 177, 3, 0, 2),
    @(54, 'Egipto', 327, 0, 56, 257, 0, 0, 14),
    @(55, 'Mexico', 316, 65, 4, 310, 1, 0, 2),
    @(56, 'Panama', 313, 0, 1, 309, 7, 0, 3),
    @(57, 'Croacia', 306, 52, 5, 300, 5, 0, 1),
    @(58, 'Libano', 267, 19, 8, 255, 4, 0, 4),
    @(59, 'Argentina', 266, 0, 27, 235, 0, 0, 4),
    @(60, 'Irak', 266, 33, 62, 181, 0, 3, 23),
    @(61, 'Colombia', 235, 4, 3, 230, 0, 0, 2),
    @(62, 'Serbia', 222, 0, 2, 218, 4, 0, 2),
    @(63, 'Republica Dominicana', 202, 0, 0, 199, 0, 0, 3),
    @(64, 'Argelia', 201, 0, 65, 119, 0, 0, 17),
    @(65, 'Taiwan', 195, 26, 28, 165, 0, 0, 2),
    @(66, 'Armenia', 194, 0, 2, 192, 6, 0, 0),
    @(67, 'Bulgaria', 190, 3, 3, 184, 8, 0, 3),
    @(68, 'Kuwait', 189, 1, 30, 159, 5, 0, 0),
    @(69, 'San Marino', 187, 12, 4, 163, 13, 0, 20),
    @(70, 'Eslovaquia', 186, 1, 7, 179, 2, 0, 0),
    @(71, 'Letonia', 180, 41, 1, 179, 0, 0, 0),
    @(72, 'Hungria', 167, 36, 16, 144, 6, 1, 7),
    @(73, 'Uruguay', 158, 0, 0, 158, 2, 0, 0),
    @(74, 'Lituania', 154, 11, 1, 152, 1, 0, 1),
    @(75, 'Emiratos Arabes Unidos', 153, 0, 38, 113, 2, 0, 2),
    @(76, 'Republica de Macedonia', 136, 21, 1, 133, 1, 1, 2),
    @(77, 'Costa Rica', 134, 0, 2, 130, 2, 0, 2),
    @(78, 'Marruecos', 134, 19, 3, 127, 1, 0, 4),
    @(79, 'Bosnia y Herzegovina', 128, 2, 2, 125, 1, 0, 1),
    @(80, 'Vietnam', 122, 9, 17, 105, 2, 0, 0),
    @(81, 'Islas Feroe', 118, 3, 14, 104, 0, 0, 0),
    @(82, 'Principado de Andorra', 113, 0, 1, 111, 2, 0, 1),
    @(83, 'Jordania', 112, 0, 1, 111, 0, 0, 0),
    @(84, 'Malta', 107, 17, 2, 105, 1, 0, 0),
    @(85, 'Nueva Zelanda', 102, 0, 0, 102, 0, 0, 0),
    @(86, 'Burkina Faso', 99, 24, 5, 90, 0, 0, 4),
    @(87, 'Republica de Chipre', 95, 0, 3, 91, 3, 0, 1),
    @(88, 'Moldavia', 94, 0, 2, 91, 3, 0, 1),
    @(89, 'Sri Lanka', 92, 10, 3, 89, 2, 0, 0),
    @(90, 'Brunei', 91, 3, 2, 89, 2, 0, 0),
    @(91, 'Tunez', 89, 14, 1, 85, 11, 0, 3),
    @(92, 'Albania', 89, 0, 2, 83, 2, 0, 0),
    @(93, 'Camboya', 86, 2, 2, 84, 1, 0, 0),
    @(94, 'Bielorrusia', 81, 5, 22, 59, 0, 0, 0),
    @(95, 'Venezuela', 77, 7, 15, 62, 2, 0, 0),
    @(96, 'Ucrania', 73, 0, 1, 69, 0, 0, 3),
    @(97, 'Azerbaiyan', 72, 7, 11, 60, 0, 0, 1),
    @(98, 'Senegal', 67, 0, 5, 62, 0, 0, 0),
    @(99, 'Oman', 66, 11, 17, 49, 0, 0, 0),
    @(100, 'Reunion', 64, 0, 1, 63, 0, 0, 0),
    @(101, 'Kazajistan', 62, 2, 0, 62, 0, 0, 0),
    @(102, 'Estado de Palestina', 59, 0, 17, 42, 0, 0, 0),
    @(103, 'Guadalupe', 58, 0, 0, 57, 4, 0, 1),
    @(104, 'Camerun', 56, 16, 2, 54, 0, 0, 0),
    @(105, 'Georgia', 54, 0, 8, 46, 1, 0, 0),
    @(106, 'Trinidad yTobago', 50, 0, 0, 50, 0, 0, 0),
    @(107, 'Liechtenstein', 46, 9, 0, 46, 0, 0, 0),
    @(108, 'Uzbekistan', 46, 3, 0, 46, 0, 0, 0),
    @(109, 'Martinica', 44, 0, 0, 43, 7, 0, 1),
    @(110, 'Afganistan', 40, 0, 1, 38, 0, 0, 1),
    @(111, 'Nigeria', 36, 6, 2, 33, 0, 1, 1),
    @(112, 'Cuba', 35, 0, 0, 34, 0, 0, 1),
    @(113, 'Banglades', 33, 6, 5, 25, 0, 1, 3),
    @(114, 'Puerto Rico', 31, 8, 0, 29, 0, 1, 2),
    @(115, 'Consejo Danes para los Refugiados', 30, 0, 0, 29, 0, 0, 1),
    @(116, 'Guam', 29, 2, 0, 28, 0, 0, 1),
    @(117, 'Mauricio', 28, 0, 0, 26, 1, 0, 2),
    @(118, 'Honduras', 27, 1, 0, 27, 0, 0, 0),
    @(119, 'Bolivia', 27, 3, 0, 27, 0, 0, 0),
    @(120, 'Costa de Marfil', 25, 11, 2, 23, 0, 0, 0),
    @(121, 'Ghana', 24, 1, 0, 23, 0, 0, 1),
    @(122, 'Macao', 24, 2, 10, 14, 0, 0, 0),
    @(123, 'Monaco', 23, 0, 1, 22, 0, 0, 0),
    @(124, 'Montenegro', 22, 1, 0, 21, 0, 1, 1),
    @(125, 'Paraguay', 22, 0, 0, 21, 1, 0, 1),
    @(126, 'Mayotte', 21, 10, 0, 21, 0, 0, 0),
    @(127, 'Ruanda', 19, 0, 0, 19, 0, 0, 0),
    @(128, 'Guyana', 19, 0, 0, 18, 0, 0, 1),
    @(129, 'Guatemala', 19, 0, 0, 18, 0, 0, 1),
    @(130, 'Jamaica', 19, 0, 2, 16, 0, 0, 1),
    @(131, 'Polinesia Francesa', 18, 0, 0, 18, 0, 0, 0),
    @(132, 'Guayana Francesa', 18, 0, 0, 18, 0, 0, 0),
    @(133, 'Togo', 18, 2, 0, 18, 0, 0, 0),
    @(134, 'Barbados', 17, 3, 0, 17, 0, 0, 0),
    @(135, 'Kirguistan', 16, 2, 0, 16, 0, 0, 0),
    @(136, 'Kenia', 15, 0, 0, 15, 0, 0, 0),
    @(137, 'Gibraltar', 15, 0, 2, 13, 0, 0, 0),
    @(138, 'Maldivas', 13, 0, 5, 8, 0, 0, 0),
    @(139, 'Madagascar', 12, 9, 0, 12, 0, 0, 0),
    @(140, 'Tanzania', 12, 0, 0, 12, 0, 0, 0),
    @(141, 'Etiopia', 11, 0, 0, 11, 0, 0, 0),
    @(142, 'Mongolia', 10, 0, 0, 10, 0, 0, 0),
    @(143, 'Aruba', 9, 0, 1, 8, 0, 0, 0),
    @(144, 'Nueva Caledonia', 8, 4, 0, 8, 0, 0, 0),
    @(145, 'Seychelles', 7, 0, 0, 7, 0, 0, 0),
    @(146, 'Guinea Ecuatorial', 6, 0, 0, 6, 0, 0, 0),
    @(147, 'Islas Virgenes de los Estados Unidos', 6, 0, 0, 6, 0, 0, 0),
    @(148, 'Bermudas', 6, 0, 0, 6, 0, 0, 0),
    @(149, 'Isla de Man', 5, 0, 0, 5, 0, 0, 0),
    @(150, 'Haiti', 5, 3, 0, 5, 0, 0, 0),
    @(151, 'Surinam', 5, 0, 0, 5, 0, 0, 0),
    @(152, 'San Martin (Parte Francesa)', 5, 0, 0, 5, 0, 0, 0),
    @(153, 'Gabon', 5, 0, 0, 4, 0, 0, 1),
    @(154, 'Bahamas', 4, 0, 0, 4, 0, 0, 0),
    @(155, 'Groenlandia', 4, 0, 0, 4, 0, 0, 0),
    @(156, 'Suazilandia', 4, 0, 0, 4, 0, 0, 0),
    @(157, 'Guinea', 4, 2, 0, 4, 0, 0, 0),
    @(158, 'Curazao', 4, 1, 0, 3, 0, 0, 1),
    @(159, 'El Salvador', 3, 0, 0, 3, 0, 0, 0),
    @(160, 'Republica de Africa Central', 3, 0, 0, 3, 0, 0, 0),
    @(161, 'Zambia', 3, 0, 0, 3, 0, 0, 0),
    @(162, 'Liberia', 3, 0, 0, 3, 0, 0, 0),
    @(163, 'Congo', 3, 0, 0, 3, 0, 0, 0),
    @(164, 'Fiyi', 3, 1, 0, 3, 0, 0, 0),
    @(165, 'San Bartolome', 3, 0, 0, 3, 0, 0, 0),
    @(166, 'Namibia', 3, 0, 0, 3, 0, 0, 0),
    @(167, 'Cabo Verde', 3, 0, 0, 3, 0, 0, 0),
    @(168, 'Islas Caimanes', 3, 0, 0, 2, 0, 0, 1),
    @(169, 'Zimbabue', 3, 0, 0, 2, 0, 1, 1),
    @(170, 'Benin', 2, 0, 0, 2, 0, 0, 0),
    @(171, 'Mauritania', 2, 0, 0, 2, 0, 0, 0),
    @(172, 'Niger', 2, 0, 0, 2, 0, 0, 0),
    @(173, 'Angola', 2, 0, 0, 2, 0, 0, 0),
    @(174, 'Santa Lucia', 2, 0, 0, 2, 0, 0, 0),
    @(175, 'Butan', 2, 0, 0, 2, 0, 0, 0),
    @(176, 'Nicaragua', 2, 0, 0, 2, 0, 0, 0),
    @(177, 'Nepal', 2, 1, 1, 1, 0, 0, 0),
    @(178, 'Gambia', 2, 1, 0, 1, 0, 1, 1),
    @(179, 'Sudan', 2, 0, 0, 1, 0, 0, 1),
    @(180, 'San Vicente y las Granadinas', 1, 0, 0, 1, 0, 0, 0),
    @(181, 'Siria', 1, 0, 0, 1, 0, 0, 0),
    @(182, 'Papua Nueva Guinea', 1, 0, 0, 1, 0, 0, 0),
    @(183, 'San Martin (Parte Holandesa)', 1, 0, 0, 1, 0, 0, 0),
    @(184, 'Timor Oriental', 1, 0, 0, 1, 0, 0, 0),
    @(185, 'Republica de Yibuti', 1, 0, 0, 1, 0, 0, 0),
    @(186, 'Eritrea', 1, 0, 0, 1, 0, 0, 0),
    @(187, 'Uganda', 1, 0, 0, 1, 0, 0, 0),
    @(188, 'Montserrat', 1, 0, 0, 1, 0, 0, 0),
    @(189, 'Mozambique', 1, 0, 0, 1, 0, 0, 0),
    @(190, 'Republica del Chad', 1, 0, 0, 1, 0, 0, 0),
    @(191, 'Somalia', 1, 0, 0, 1, 0, 0, 0),
    @(192, 'Antigua y Barbuda', 1, 0, 0, 1, 0, 0, 0),
    @(193, 'Dominica', 1, 0, 0, 1, 0, 0, 0),
    @(194, 'Santa Sede', 1, 0, 0, 1, 0, 0, 0)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
}
